# "signature" -> "Signature" in the table header cell, but Word records
# the edit as two adjacent runs (the capitalized "S" typed in place of the
# old "s", followed by the untouched "ignature") rather than merging back
# into a single run, even though both runs end up with identical
# formatting.
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("signature", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $first = $d.Range($rng.Start, $rng.Start + 1)

    # Replace the lowercase "s" with an uppercase "S".
    $first.Text = "S"

    # Toggling bold off/on for just this new character forces Word to
    # keep it in its own run instead of silently re-merging it with the
    # following "ignature" run, matching the split seen in the diff.
    $first2 = $d.Range($rng.Start, $rng.Start + 1)
    $first2.Font.Bold = 0
    $first2.Font.Bold = 1
}
